$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "WHEAT" title shape
$wheatShape = $s.Shapes.Item(1)
$wheatFont = $wheatShape.TextFrame.TextRange.Font
$wheatFont.Bold = $true
$wheatFont.Name = "+mn-lt"

# "BLACK" title shape
$blackShape = $s.Shapes.Item(4)
$blackFont = $blackShape.TextFrame.TextRange.Font
$blackFont.Bold = $true
$blackFont.Name = "+mn-lt"
